# Updates odds values on Sheet1 (Jogos da Semana FlashScore 2024-11-08)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (match 1)
$ws.Range("G2").Value = 1.55
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("AC2").Value = 8
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81

# Row 8 (match 7)
$ws.Range("G8").Value = 1.14
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 13
$ws.Range("W8").Value = 11
$ws.Range("AB8").Value = 29
$ws.Range("AC8").Value = 23
$ws.Range("AE8").Value = 29
$ws.Range("AG8").Value = 1000
$ws.Range("AL8").Value = 67
$ws.Range("AM8").Value = 51
$ws.Range("AQ8").Value = 10

# Row 9 (match 8)
$ws.Range("G9").Value = 1.22
$ws.Range("H9").Value = 6.5
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 1.62
$ws.Range("L9").Value = 8
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 13
$ws.Range("Q9").Value = 1.4
$ws.Range("R9").Value = 2.75
$ws.Range("S9").Value = 1.22
$ws.Range("T9").Value = 4
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 10
$ws.Range("Y9").Value = 10
$ws.Range("Z9").Value = 8
$ws.Range("AB9").Value = 26
$ws.Range("AG9").Value = 800
$ws.Range("AH9").Value = 26
$ws.Range("AJ9").Value = 26
$ws.Range("AQ9").Value = 12
$ws.Range("AT9").Value = 4
$ws.Range("AW9").Value = 11

$wb.Save()
